$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "2025-09-05" parses as a real date if written directly, which would
# stamp the cell with a date number format. Force text ("@") first so the
# literal string is preserved, then restore the default "Normal" style so
# the new row doesn't end up with a leftover explicit style index (keeps
# it consistent with the rest of the sheet, which has no per-cell styles).
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2025-09-05"
$ws.Range("A6").Style = "Normal"

$ws.Range("B6").Value = "21:21:03"
$ws.Range("C6").Value = "1.00 EUR = 1611.0529 ARS"
